{"js": "// Apply the text edits described by the commit:\n// \"Deleted Escape package from Italy\" \u2014 a copy-edit pass over the\n// interview transcript paragraphs. We locate each changed phrase with\n// body.search() (exact text, case sensitive) and replace it in place.\n\nasync function replaceOnce(context, findText, replaceText) {\n  const results = context.document.body.search(findText, {\n    matchCase: true,\n    matchWholeWord: false,\n  });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + findText);\n  }\n\n  // Only the first occurrence is expected/needed for each of these edits.\n  results.items[0].insertText(replaceText, \"Replace\");\n  await context.sync();\n}\n\n// 1) \"Clayton:  We started with sketching out our storyboard...\" paragraph\nawait replaceOnce(\n  context,\n  \"We decided we\\u2019d made a slider\",\n  \"We decided we\\u2019d make a slider\"\n);\n\nawait replaceOnce(\n  context,\n  \"We then decided the countries we would allow customers to visit, deciding on Italy, Japan, and Jamaica.\",\n  \"We then decided that the countries would be Italy, Japan, and Jamaica.\"\n);\n\n// 2) \"Tyler:  We then began work creating the actual page...\" paragraph\nawait replaceOnce(\n  context,\n  \"Clayton began working on the slider using Jquery.  \",\n  \"\"\n);\n\n// 3) \"Tyler:  I then started finding images for Italy...\" paragraph\nawait replaceOnce(\n  context,\n  \"We, as a business would have partnered\",\n  \"We, as a business have partnered\"\n);\n\nawait replaceOnce(\n  context,\n  \"Finally the rent-a-car company I chose would be Enterprise, I chose them because\",\n  \"Finally the rent-a-car company would be Enterprise. I selected enterprise because\"\n);\n\n// 4) \"Tyler:  We then began to add finishing touches on the site...\" paragraph\nawait replaceOnce(context, \"I wrote scripts\", \"I then wrote text scripts\");\n\n// 5) \"Tyler:  We then began to work on making sure...\" paragraph (last one)\nawait replaceOnce(\n  context,\n  \"we hit everything on the rubric.  Cross browser\",\n  \"the website had  cross browser\"\n);\n\nawait replaceOnce(context, \"validating etc.\", \"validating the code.\");\n\n// The last edit position (the \"_GoBack\" bookmark Word maintains\n// automatically) moves from the \"destinations page\" paragraph to the\n// very end of the document, since that's where the final edit above\n// took place.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst lastParagraph = context.document.body.paragraphs.getLast();\nconst endOfDoc = lastParagraph.getRange(\"End\");\nendOfDoc.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "function Replace-Text {\n    param($doc, [string]$findText, [string]$replaceText)\n\n    $find = $doc.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $replaceText\n    # wdFindContinue=1, wdReplaceAll=2 (each search phrase below is unique\n    # in the document, so \"replace all\" behaves the same as \"replace one\").\n    $found = $find.Execute($findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)\n    if (-not $find.Found) {\n        throw \"Text not found: $findText\"\n    }\n}\n\n$d = $word.ActiveDocument\n\n# 1) \"Clayton:  We started with sketching out our storyboard...\" paragraph\nReplace-Text $d \"We decided we\u2019d made a slider\" \"We decided we\u2019d make a slider\"\nReplace-Text $d \"We then decided the countries we would allow customers to visit, deciding on Italy, Japan, and Jamaica.\" \"We then decided that the countries would be Italy, Japan, and Jamaica.\"\n\n# 2) \"Tyler:  We then began work creating the actual page...\" paragraph\nReplace-Text $d \"Clayton began working on the slider using Jquery.  \" \"\"\n\n# 3) \"Tyler:  I then started finding images for Italy...\" paragraph\nReplace-Text $d \"We, as a business would have partnered\" \"We, as a business have partnered\"\nReplace-Text $d \"Finally the rent-a-car company I chose would be Enterprise, I chose them because\" \"Finally the rent-a-car company would be Enterprise. I selected enterprise because\"\n\n# 4) \"Tyler:  We then began to add finishing touches on the site...\" paragraph\nReplace-Text $d \"I wrote scripts\" \"I then wrote text scripts\"\n\n# 5) \"Tyler:  We then began to work on making sure...\" paragraph (last one)\nReplace-Text $d \"we hit everything on the rubric.  Cross browser\" \"the website had  cross browser\"\nReplace-Text $d \"validating etc.\" \"validating the code.\"\n\nWrite-Output \"done\"\n"}
